$d = $word.ActiveDocument

# --- Change 1: fix "following: -" to "following: " and remove the now-stray empty paragraph after it ---
$findRange = $d.Content
$found = $findRange.Find.Execute("following: -", $false, $false, $false, $false, $false, $true, 1, $false, "following: ", 2)
if ($found) {
    $fixedPara = $findRange.Paragraphs(1)
    $nextPara = $fixedPara.Next()
    if ($nextPara -ne $null -and $nextPara.Range.Text.Trim().Length -eq 0) {
        $nextPara.Range.Delete()
    }
}

# --- Change 2: append the new "Removing custom profile" / "Loading profile" sections at the end of the body ---
$bodyFragment = '<w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Removing custom profile in Obi:</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">To remove </w:t></w:r><w:r><w:t>profile in combo box of Select profile follow following steps:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Step 1: </w:t></w:r><w:r><w:t xml:space="preserve">Select profile from combo box you want to </w:t></w:r><w:r><w:t>delete</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>Step 2: Now click remove button.</w:t></w:r></w:p><w:p><w:r><w:t>Step 3: A warning message stating “</w:t></w:r><w:r><w:t>This will delete the selected profile. Do you want to continue?</w:t></w:r><w:r><w:t>” will come. Press yes to continue.</w:t></w:r></w:p><w:p><w:r><w:t>Step 4: Selected profile will be deleted.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>Loading profile in Obi:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:r><w:t xml:space="preserve">To load </w:t></w:r><w:r><w:t>profile in</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Obi</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> follow following steps:</w:t></w:r></w:p><w:p><w:r><w:t>Step 1:</w:t></w:r><w:r><w:t xml:space="preserve"> Select profile from combo box you want to load. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Step 2: On the right side of select profile combo box you can find list of check boxes for </w:t></w:r><w:r><w:t>Project / Audio / Language / Color / Select All</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Check the check boxes for the preferences that you wish to load.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Step 3: </w:t></w:r><w:r><w:t xml:space="preserve">Click on load profile button. This will load selected profile on Obi. </w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Note</w:t></w:r><w:r><w:t>: I</w:t></w:r><w:r><w:t>f the check box pertinent to the set of preferences is not checked, they will not be loaded; for example, if the Project check box is not checked, the project preferences will not be loaded.</w:t></w:r></w:p><w:p/><w:p/><w:p/>'

$xmlPayload = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($xmlPayload)
